# Add a new "2022-Q3" quarter sheet, shift the summary table, and
# insert fund-holding detail rows for the new quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: insert the 2022-Q3 row at
#    the top of the data, and append the 2021-Q2 row that used to
#    be the last row but now needs an extra row because everything
#    else shifted down by one.
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 13
$summary.Cells.Item(2, 4).Value = 0.63

$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(3, 2).Value = "2022-Q2"
$summary.Cells.Item(3, 3).Value = 8
$summary.Cells.Item(3, 4).Value = 0.13

$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(4, 2).Value = "2022-Q1"
$summary.Cells.Item(4, 3).Value = 2
$summary.Cells.Item(4, 4).Value = 0.23

$summary.Cells.Item(5, 1).Value = 3
$summary.Cells.Item(5, 2).Value = "2021-Q3"
$summary.Cells.Item(5, 3).Value = 1
$summary.Cells.Item(5, 4).Value = 0.45

$summary.Cells.Item(6, 1).Value = 4
$summary.Cells.Item(6, 2).Value = "2021-Q2"
$summary.Cells.Item(6, 3).Value = 1
$summary.Cells.Item(6, 4).Value = 0.38

# Column A on the new row needs the same style as the rest of column A
# (bold / centered / bordered look used throughout the workbook).
# Copy the format down from the row above it.
$summary.Range("A5").Copy()
$summary.Range("A6").PasteSpecial(-4122)

# ---------------------------------------------------------------
# 2. Insert the brand-new "2022-Q3" worksheet right after "总计"
#    (i.e. before "2022-Q2", which is currently the 2nd sheet).
#
#    A freshly-added blank worksheet doesn't reliably pick up
#    pasted formatting in this runtime, so instead duplicate the
#    existing, already-styled "2022-Q2" sheet (which shares the
#    exact same header/column layout) and overwrite its values.
# ---------------------------------------------------------------
$q2sheet = $wb.Worksheets.Item("2022-Q2")
$q2sheet.Copy($null, $summary)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# The duplicated sheet only has 8 data rows (rows 2-9); we need 13
# (rows 2-14). Extend column A's styling down to the new rows.
$newSheet.Range("A9").Copy()
$newSheet.Range("A10:A14").PasteSpecial(-4122)

function Set-FundRow {
    param($sheet, $row, $idx, $code, $fundName, $size, $pos, $ratio, $value, $rank, $valueIsText)

    $sheet.Cells.Item($row, 1).Value = $idx

    # Fund code, scale, position-size and position-ratio columns are
    # stored as plain text in the source data (so things like leading
    # zeros in the fund code survive) - force text before assigning.
    $sheet.Cells.Item($row, 2).NumberFormat = "@"
    $sheet.Cells.Item($row, 2).Value = $code
    $sheet.Cells.Item($row, 3).Value = $fundName
    $sheet.Cells.Item($row, 4).NumberFormat = "@"
    $sheet.Cells.Item($row, 4).Value = $size
    $sheet.Cells.Item($row, 5).NumberFormat = "@"
    $sheet.Cells.Item($row, 5).Value = $pos
    $sheet.Cells.Item($row, 6).NumberFormat = "@"
    $sheet.Cells.Item($row, 6).Value = $ratio
    if ($valueIsText) {
        $sheet.Cells.Item($row, 7).NumberFormat = "@"
        $sheet.Cells.Item($row, 7).Value = $value
    } else {
        $sheet.Cells.Item($row, 7).Value = [double]$value
    }
    $sheet.Cells.Item($row, 8).Value = $rank
}

Set-FundRow $newSheet 2  0 "001468" "广发改革先锋灵活配置混合"     "6.25" "62.12" "4.90" "0.3062" 2 $true
Set-FundRow $newSheet 3  1 "090016" "大成消费主题混合"             "3.59" "90.85" "4.44" "0.1594" 8 $true
Set-FundRow $newSheet 4  2 "860058" "光大阳光稳健增长混合C"         "6.74" "32.18" "0.65" "0.0438" 8 $true
Set-FundRow $newSheet 5  3 "860009" "光大阳光稳健增长混合A"         "6.09" "32.18" "0.65" "0.0396" 8 $true
Set-FundRow $newSheet 6  4 "013204" "恒生前海恒源天利债A"          "1.30" "32.00" "1.80" "0.0234" 3 $true
Set-FundRow $newSheet 7  5 "860006" "光大阳光优选一年持有混合A"     "1.04" "79.87" "2.10" "0.0218" 7 $true
Set-FundRow $newSheet 8  6 "014151" "国富鑫享价值一年封闭混合A"     "2.16" "46.24" "0.88" "0.0190" 9 $true
Set-FundRow $newSheet 9  7 "011590" "九泰天利量化股票C"             "0.50" "83.77" "2.26" "0.0113" 3 $true
Set-FundRow $newSheet 10 8 "014152" "国富鑫享价值一年封闭混合C"     "0.75" "46.24" "0.88" "0.0066" 9 $true
Set-FundRow $newSheet 11 9 "011589" "九泰天利量化股票A"             "0.07" "83.77" "2.26" "0.0016" 3 $true
Set-FundRow $newSheet 12 10 "860055" "光大阳光优选一年持有混合B"    "0.00" "79.87" "2.10" 0        7 $false
Set-FundRow $newSheet 13 11 "860056" "光大阳光优选一年持有混合C"    "0.00" "79.87" "2.10" 0        7 $false
Set-FundRow $newSheet 14 12 "013205" "恒生前海恒源天利债C"         "0.00" "32.00" "1.80" 0        3 $false

$newSheet.Range("A1").Select()
